$wb = $excel.ActiveWorkbook

# Cosmetic: the sheet-tab/horizontal-scrollbar splitter ratio moved slightly
# in the authored workbook (973 -> 992 on a 0-1000 scale in the OOXML, i.e.
# ~0.973 -> ~0.992 via the COM TabRatio property).
try { $excel.ActiveWindow.TabRatio = 0.992 } catch { }

# --- Commodity sheet: add a 4th row (Left / Magic / 5) ---
$wsCommodity = $wb.Worksheets.Item("Commodity")
$wsCommodity.Range("A4").Value = "Left"
$wsCommodity.Range("B4").Value = "Magic"
$wsCommodity.Range("C4").Value = 5
$wsCommodity.Range("C1").Value = "Price"

# --- Process sheet: add MinOut/MaxOut columns and a new "Mage" process row ---
$wsProcess = $wb.Worksheets.Item("Process")
$wsProcess.Range("C1").Value = "MinOut"
$wsProcess.Range("D1").Value = "MaxOut"

$wsProcess.Range("A2").Value = "Left"
$wsProcess.Range("B2").Value = "Coal plant"
$wsProcess.Range("C2").Value = 0
$wsProcess.Range("D2").Value = 1000

$wsProcess.Range("A3").Value = "Left"
$wsProcess.Range("B3").Value = "Mage"
$wsProcess.Range("C3").Value = 0
$wsProcess.Range("D3").Value = 5

$wsProcess.Range("A4").Value = "Right"
$wsProcess.Range("B4").Value = "Coal plant"
$wsProcess.Range("C4").Value = 0
$wsProcess.Range("D4").Value = 1000

# --- Process-Commodity sheet: add Mage in/out rows ---
$wsPC = $wb.Worksheets.Item("Process-Commodity")
$wsPC.Range("A4").Value = "Mage"
$wsPC.Range("B4").Value = "Magic"
$wsPC.Range("C4").Value = "in"
$wsPC.Range("D4").Value = 1

$wsPC.Range("A5").Value = "Mage"
$wsPC.Range("B5").Value = "Elec"
$wsPC.Range("C5").Value = "out"
$wsPC.Range("D5").Value = 1

# --- Restore the originally active sheet/selection state ---
# Process sheet selection ends on C2 (per target), Process-Commodity ends on B5 (per target),
# but the workbook's active/visible tab must remain "Commodity" at C1 (unchanged from before).
$wsProcess.Range("C2").Select() | Out-Null
$wsPC.Range("B5").Select() | Out-Null
$wsCommodity.Activate() | Out-Null
$wsCommodity.Range("C1").Select() | Out-Null

$wb.Save()
